$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 33354
$ws.Range("J87").Value = 33354
$ws.Range("L87").Value = 33354
$ws.Range("N87").Value = -35850
$ws.Range("H90").Value = 33354
$ws.Range("J90").Value = 33354
$ws.Range("L90").Value = 100062
$ws.Range("N90").Value = -112542
$ws.Range("H106").Value = 10472.583
$ws.Range("I106").Value = 1959.625
$ws.Range("K106").Value = 1959.625
$ws.Range("M106").Value = -1328.625
$ws.Range("H116").Value = 4999.8
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 4999.75
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 4999.75
$ws.Range("M116").Value = -1558
$ws.Range("N116").Value = -11883.75
$ws.Range("H131").Value = 750
$ws.Range("I131").Value = 750
$ws.Range("K131").Value = 2250
$ws.Range("M131").Value = 2790
$ws.Range("H132").Value = 1281.7021
$ws.Range("I132").Value = 868.3953
$ws.Range("J132").Value = 5724.75
$ws.Range("K132").Value = 2605.1859
$ws.Range("L132").Value = 17174.25
$ws.Range("M132").Value = -75.18589999999995
$ws.Range("N132").Value = -22234.25
$ws.Range("H138").Value = 2456.5676
$ws.Range("I138").Value = 1480.7142
$ws.Range("J138").Value = 3050.5652
$ws.Range("K138").Value = 4442.142599999999
$ws.Range("L138").Value = 9151.695599999999
$ws.Range("M138").Value = 697.8574000000008
$ws.Range("N138").Value = -19431.6956

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1241.027
$ws.Range("I61").Value = 1136.6111
$ws.Range("K61").Value = 1136.6111
$ws.Range("M61").Value = -924.6111000000001
$ws.Range("H88").Value = 808
$ws.Range("I88").Value = 475
$ws.Range("J88").Value = 941.2
$ws.Range("K88").Value = 475
$ws.Range("L88").Value = 941.2
$ws.Range("M88").Value = -69
$ws.Range("N88").Value = -1753.2
$ws.Range("H91").Value = 808
$ws.Range("I91").Value = 475
$ws.Range("J91").Value = 941.2
$ws.Range("K91").Value = 475
$ws.Range("L91").Value = 941.2
$ws.Range("M91").Value = 929
$ws.Range("N91").Value = -3749.2
$ws.Range("H122").Value = 1773.4286
$ws.Range("I122").Value = 1637.1
$ws.Range("K122").Value = 4911.299999999999
$ws.Range("M122").Value = -2461.299999999999
$ws.Range("H132").Value = 2103.4211
$ws.Range("I132").Value = 1186.5312
$ws.Range("K132").Value = 3559.5936
$ws.Range("M132").Value = -1029.5936
$ws.Range("H135").Value = 63499.5
$ws.Range("J135").Value = 63499.5
$ws.Range("L135").Value = 63499.5
$ws.Range("N135").Value = -73639.5
$ws.Range("H136").Value = 1241.027
$ws.Range("I136").Value = 1136.6111
$ws.Range("K136").Value = 3409.8333
$ws.Range("M136").Value = -859.8333000000002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 33332.332
$ws.Range("J81").Value = 33332.332
$ws.Range("L81").Value = 33332.332
$ws.Range("N81").Value = -35454.332
$ws.Range("H82").Value = 17917.445
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766
$ws.Range("H84").Value = 33332.332
$ws.Range("J84").Value = 33332.332
$ws.Range("L84").Value = 99996.99600000001
$ws.Range("N84").Value = -110604.996
$ws.Range("H85").Value = 17917.445
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652
$ws.Range("H94").Value = 1263.5
$ws.Range("I94").Value = 1299.7142
$ws.Range("K94").Value = 1299.7142
$ws.Range("M94").Value = -848.7141999999999
$ws.Range("H102").Value = 11221
$ws.Range("I102").Value = 1525.5
$ws.Range("J102").Value = 30612
$ws.Range("K102").Value = 1525.5
$ws.Range("L102").Value = 30612
$ws.Range("M102").Value = 1719.5
$ws.Range("N102").Value = -37102
$ws.Range("H107").Value = 1146.3125
$ws.Range("I107").Value = 1205.0769
$ws.Range("K107").Value = 1205.0769
$ws.Range("M107").Value = 714.9231
$ws.Range("H134").Value = 1469.6786
$ws.Range("I134").Value = 1450.037
$ws.Range("K134").Value = 4350.111
$ws.Range("M134").Value = -1815.111

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4350.6665
$ws.Range("I31").Value = 1926.5
$ws.Range("J31").Value = 7380.875
$ws.Range("K31").Value = 1926.5
$ws.Range("L31").Value = 7380.875
$ws.Range("M31").Value = -1631.5
$ws.Range("N31").Value = -7970.875
$ws.Range("H34").Value = 4350.6665
$ws.Range("I34").Value = 1926.5
$ws.Range("J34").Value = 7380.875
$ws.Range("K34").Value = 1926.5
$ws.Range("L34").Value = 7380.875
$ws.Range("M34").Value = -1724.5
$ws.Range("N34").Value = -7784.875
$ws.Range("H58").Value = 1373.1666
$ws.Range("I58").Value = 1395.4375
$ws.Range("J58").Value = 1328.625
$ws.Range("K58").Value = 1395.4375
$ws.Range("L58").Value = 1328.625
$ws.Range("M58").Value = -1192.4375
$ws.Range("N58").Value = -1734.625
$ws.Range("H68").Value = 25000
$ws.Range("J68").Value = 25000
$ws.Range("L68").Value = 25000
$ws.Range("N68").Value = -26498
$ws.Range("H71").Value = 25000
$ws.Range("J71").Value = 25000
$ws.Range("L71").Value = 75000
$ws.Range("N71").Value = -82488
$ws.Range("H74").Value = 38701.668
$ws.Range("J74").Value = 56663
$ws.Range("L74").Value = 56663
$ws.Range("N74").Value = -58411
$ws.Range("H77").Value = 38701.668
$ws.Range("J77").Value = 56663
$ws.Range("L77").Value = 169989
$ws.Range("N77").Value = -178725
$ws.Range("H95").Value = 13349.667
$ws.Range("J95").Value = 13349.667
$ws.Range("L95").Value = 13349.667
$ws.Range("N95").Value = -18841.667
$ws.Range("H134").Value = 1897.2
$ws.Range("I134").Value = 1824.6897
$ws.Range("K134").Value = 5474.0691
$ws.Range("M134").Value = -2939.0691
$ws.Range("H136").Value = 1373.1666
$ws.Range("I136").Value = 1395.4375
$ws.Range("J136").Value = 1328.625
$ws.Range("K136").Value = 4186.3125
$ws.Range("L136").Value = 3985.875
$ws.Range("M136").Value = -1636.3125
$ws.Range("N136").Value = -9085.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 10333.333
$ws.Range("J105").Value = 15000
$ws.Range("L105").Value = 45000
$ws.Range("N105").Value = -50242
$ws.Range("H138").Value = 1000
$ws.Range("I138").Value = 1000
$ws.Range("K138").Value = 3000
$ws.Range("M138").Value = 2140
$ws.Range("H140").Value = 3710.5386
$ws.Range("I140").Value = 2279.875
$ws.Range("K140").Value = 6839.625
$ws.Range("M140").Value = -1659.625
$ws.Range("H141").Value = 3099.8572
$ws.Range("I141").Value = 3116.6667
$ws.Range("K141").Value = 9350.000100000001
$ws.Range("M141").Value = -4170.000100000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 24496
$ws.Range("I36").Value = 24496
$ws.Range("K36").Value = 24496
$ws.Range("M36").Value = -24011
$ws.Range("H126").Value = 6817.231
$ws.Range("I126").Value = 9135.429
$ws.Range("J126").Value = 4112.6665
$ws.Range("K126").Value = 27406.287
$ws.Range("L126").Value = 12337.9995
$ws.Range("M126").Value = -24936.287
$ws.Range("N126").Value = -17277.9995
$ws.Range("H132").Value = 1934.6451
$ws.Range("I132").Value = 1771.76
$ws.Range("J132").Value = 2613.3333
$ws.Range("K132").Value = 5315.28
$ws.Range("L132").Value = 7839.999899999999
$ws.Range("M132").Value = -2785.28
$ws.Range("N132").Value = -12899.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 87.05882
$ws.Range("J2").Value = 87.05882
$ws.Range("L2").Value = 87.05882
$ws.Range("N2").Value = -311.05882
$ws.Range("H40").Value = 3522.182
$ws.Range("I40").Value = 3074.9
$ws.Range("J40").Value = 7995
$ws.Range("K40").Value = 3074.9
$ws.Range("L40").Value = 7995
$ws.Range("M40").Value = -2938.9
$ws.Range("N40").Value = -8267
$ws.Range("H64").Value = 14691
$ws.Range("J64").Value = 14691
$ws.Range("L64").Value = 14691
$ws.Range("N64").Value = -15141
$ws.Range("H67").Value = 14691
$ws.Range("J67").Value = 14691
$ws.Range("L67").Value = 14691
$ws.Range("N67").Value = -16251
$ws.Range("H132").Value = 3682.7856
$ws.Range("I132").Value = 2784.611
$ws.Range("K132").Value = 8353.832999999999
$ws.Range("M132").Value = -5823.832999999999
$ws.Range("H136").Value = 2847.6924
$ws.Range("I136").Value = 2418.75
$ws.Range("K136").Value = 7256.25
$ws.Range("M136").Value = -4706.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 45034.5
$ws.Range("H63").Value = 59949
$ws.Range("J63").Value = 59949
$ws.Range("L63").Value = 59949
$ws.Range("N63").Value = -61197
$ws.Range("H66").Value = 59949
$ws.Range("J66").Value = 59949
$ws.Range("L66").Value = 179847
$ws.Range("N66").Value = -186087
$ws.Range("H81").Value = 4047.68
$ws.Range("I81").Value = 4047.68
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 8095.36
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -7034.36
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 4047.68
$ws.Range("I84").Value = 4047.68
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 40476.8
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -35172.8
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 4794.206
$ws.Range("I132").Value = 4919.636
$ws.Range("K132").Value = 14758.908
$ws.Range("M132").Value = -12228.908
$ws.Range("H136").Value = 200
$ws.Range("I136").Value = 200
$ws.Range("K136").Value = 600
$ws.Range("M136").Value = 1950
